$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = -129300000.0
$ws.Range("C6").Value = -69700000.0
$ws.Range("D6").Value = -43500000.0
$ws.Range("E6").Value = -59300000.0
$ws.Range("F6").Value = -37600000.0

$ws.Range("B8").Value = 2591000000.0
$ws.Range("C8").Value = 2496000000.0
$ws.Range("D8").Value = 1756100000.0
$ws.Range("E8").Value = 1191400000.0
$ws.Range("F8").Value = 548900000.0
